# Commit: "Added Type property to the Event entity (same in DB), and added
# new users statistics"
#
# 1) Add a new first sheet "New Users" with a Day/Users table (new daily
#    new-user-count report).
# 2) The "Preliminary statistics" sheet's item rows come from a DB query
#    that (because of the new Type column) now returns a handful of items
#    in a slightly different relative order. No values actually changed,
#    just which row a given item (and its Amount/Currency/Price) sits on,
#    for five small neighbouring groups of rows.

$wb = $excel.ActiveWorkbook

# --- 1. New "New Users" sheet, inserted as the very first tab ----------
$newUsers = $wb.Worksheets.Add()
$newUsers.Name = "New Users"

$newUsers.Range("A1").Value = "Day"
$newUsers.Range("B1").Value = "Users"

# Keep "01.01.2018" as text (matches the same label used on the DAU/MAU
# sheets) instead of letting it be auto-recognised as a date serial.
$newUsers.Range("A2").NumberFormat = "@"
$newUsers.Range("A2").Value = "01.01.2018"
$newUsers.Range("A2").Style = "Normal"
$newUsers.Range("B2").Value = 14831

# --- 2. Re-order a few item rows on "Preliminary statistics" -----------
$ws = $wb.Worksheets.Item("Preliminary statistics")

$ws.Range("A17").Value = "Converter"
$ws.Range("B17").Value = 14
$ws.Range("C17").Value = 9800
$ws.Range("D17").Value = 0

$ws.Range("A18").Value = "Dead Sea Scrolls"
$ws.Range("B18").Value = 11
$ws.Range("C18").Value = 5940
$ws.Range("D18").Value = 0

$ws.Range("A19").Value = "Butter Bean"
$ws.Range("B19").Value = 23
$ws.Range("C19").Value = 1380
$ws.Range("D19").Value = 0

$ws.Range("A54").Value = "The Bean"
$ws.Range("B54").Value = 16
$ws.Range("C54").Value = 1600
$ws.Range("D54").Value = 0

$ws.Range("A55").Value = "Best Friend"
$ws.Range("B55").Value = 16
$ws.Range("C55").Value = 10400
$ws.Range("D55").Value = 0

$ws.Range("A80").Value = "Mr. Boom"
$ws.Range("B80").Value = 19
$ws.Range("C80").Value = 9310
$ws.Range("D80").Value = 0

$ws.Range("A81").Value = "Void"
$ws.Range("B81").Value = 9
$ws.Range("C81").Value = 5040
$ws.Range("D81").Value = 0

$ws.Range("A82").Value = "The Boomerang"
$ws.Range("B82").Value = 11
$ws.Range("C82").Value = 7480
$ws.Range("D82").Value = 0

$ws.Range("A89").Value = "Razor Blade"
$ws.Range("B89").Value = 21
$ws.Range("C89").Value = 3570
$ws.Range("D89").Value = 0

$ws.Range("A90").Value = "Wooden Nickel"
$ws.Range("B90").Value = 11
$ws.Range("C90").Value = 1210
$ws.Range("D90").Value = 0

$ws.Range("A109").Value = "Black Hole"
$ws.Range("B109").Value = 11
$ws.Range("C109").Value = 5610
$ws.Range("D109").Value = 0

$ws.Range("A110").Value = "Sacrificial Altar"
$ws.Range("B110").Value = 7
$ws.Range("C110").Value = 4200
$ws.Range("D110").Value = 0

$ws.Range("A111").Value = "The Hourglass"
$ws.Range("B111").Value = 19
$ws.Range("C111").Value = 1520
$ws.Range("D111").Value = 0
